# Applies the scheduled data refresh:
#  1) Re-sorts the per-matchday block of rows (columns F:V only — the match
#     facts: teams, scores, odds + timestamps, source URL) while columns
#     A:E (sequence number / league / season / match date) stay put.
#  2) Appends 9 newly scraped fixtures (rows 101-109).
# Excel recomputes UsedRange / <dimension> automatically once the new cells
# are populated, so no manual dimension bookkeeping is required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Rotate-RowFacts {
    # Cycle the F:V ("match fact") payload of a set of rows forward by one
    # position: row[i] receives what used to live in row[i+1] (wrapping).
    param($ws, [int[]]$rows)

    $n = $rows.Count
    $vals = @()
    foreach ($r in $rows) {
        $vals += , ($ws.Range("F${r}:V${r}").Value2)
    }
    for ($i = 0; $i -lt $n; $i++) {
        $dst = $rows[$i]
        $src = $vals[($i + 1) % $n]
        $ws.Range("F${dst}:V${dst}").Value2 = $src
    }
}

# Groups of rows that belong to the same match day and were re-ordered.
Rotate-RowFacts $ws @(26, 28)
Rotate-RowFacts $ws @(48, 50)
Rotate-RowFacts $ws @(52, 53)
Rotate-RowFacts $ws @(57, 58, 59)
Rotate-RowFacts $ws @(70, 71)
Rotate-RowFacts $ws @(73, 75)
Rotate-RowFacts $ws @(79, 81)
Rotate-RowFacts $ws @(80, 82)
Rotate-RowFacts $ws @(89, 91)
Rotate-RowFacts $ws @(96, 97)

# --- Append newly scraped fixtures (rows 101-109) -------------------------

# Carry the row-100 formatting (bold/bordered match-number cell, date/time
# number format on column E, etc.) down onto the freshly appended rows.
$ws.Range("A100:V100").Copy()
$ws.Range("A101:V109").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$newRows = @(
    # A,   B,        C,                    D,            E,                    F,                        G, H,                      I, J,    K,                   L,    M,                   N,    O,                   P,    Q,                   R,    S,                   T,    U,                   V
    @(100, "poland", "iii-liga-group-iv", "2023-2024", 45226.79166666666, "Swidniczanka Swidnik", 2, "KS Wieczysta Krakow", 4, 5.71, "26/10/2023 07:12", 5.87, "27/10/2023 18:53", 4.9,  "26/10/2023 07:12", 5.05, "27/10/2023 18:53", 1.29, "26/10/2023 07:12", 1.36, "27/10/2023 18:53", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/swidniczanka-swidnik-ks-wieczysta-krakow/W8zrIY0p/"),
    @(101, "poland", "iii-liga-group-iv", "2023-2024", 45227.5,           "Wislanie Jaskowice",   0, "Biala Podlaska",       0, 2.27, "27/10/2023 00:12", 2.24, "28/10/2023 11:56", 3.26, "27/10/2023 00:12", 3.42, "28/10/2023 11:56", 2.5,  "27/10/2023 00:12", 2.72, "28/10/2023 11:56", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/wislanie-jaskowice-biala-podlaska/YXKtwzwA/"),
    @(102, "poland", "iii-liga-group-iv", "2023-2024", 45227.58333333334, "Chelmianka Chelm",     2, "Avia Swidnik",         2, 3.11, "27/10/2023 02:13", 3.26, "28/10/2023 13:57", 3.3,  "27/10/2023 02:13", 3.57, "28/10/2023 13:57", 1.91, "27/10/2023 02:13", 1.93, "28/10/2023 13:57", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/chelmianka-chelm-avia-swidnik/Q5xFKRMr/"),
    @(103, "poland", "iii-liga-group-iv", "2023-2024", 45227.58333333334, "Czarni Polaniec",      1, "Unia Tarnow",          2, 1.87, "27/10/2023 02:13", 1.81, "28/10/2023 07:30", 3.58, "27/10/2023 02:13", 3.85, "28/10/2023 13:23", 3,    "27/10/2023 02:13", 3.42, "28/10/2023 13:23", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/czarni-polaniec-unia-tarnow/nHIUuxgi/"),
    @(104, "poland", "iii-liga-group-iv", "2023-2024", 45227.58333333334, "Garbarnia",            1, "Wiazownica",           3, 1.9,  "27/10/2023 02:13", 1.75, "28/10/2023 07:27", 3.51, "27/10/2023 02:13", 3.84, "28/10/2023 12:01", 2.96, "27/10/2023 02:13", 3.6,  "28/10/2023 07:27", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/garbarnia-wiazownica/GdAoxfhG/"),
    @(105, "poland", "iii-liga-group-iv", "2023-2024", 45227.58333333334, "Karpaty Krosno",       3, "Orleta Radzyn",        5, 2.07, "27/10/2023 02:13", 1.99, "28/10/2023 13:49", 3.3,  "27/10/2023 02:13", 3.4,  "28/10/2023 13:49", 2.78, "27/10/2023 02:13", 3.24, "28/10/2023 13:49", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/ks-karpaty-krosno-orleta-radzyn/fsMxvGN3/"),
    @(106, "poland", "iii-liga-group-iv", "2023-2024", 45227.58333333334, "Siarka Tarnobrzeg",    3, "Wisloka Debica",       0, 1.42, "27/10/2023 02:13", 1.55, "28/10/2023 13:58", 4.12, "27/10/2023 02:13", 4.01, "28/10/2023 13:58", 4.91, "27/10/2023 02:13", 4.76, "28/10/2023 13:58", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/siarka-tarnobrzeg-wisloka-debica/pCvnHhGj/"),
    @(107, "poland", "iii-liga-group-iv", "2023-2024", 45227.60416666666, "Star Starachowice",    0, "Ostrowiec Swietokrzyski", 1, 1.98, "27/10/2023 02:42", 2.04, "28/10/2023 13:46", 3.26, "27/10/2023 02:42", 3.23, "28/10/2023 13:46", 2.99, "27/10/2023 02:42", 3.27, "28/10/2023 13:46", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/star-starachowice-ostrowiec-swietokrzyski/Q7HYvd8c/"),
    @(108, "poland", "iii-liga-group-iv", "2023-2024", 45228.58333333334, "Sokol Sieniawa",       1, "Podhale Nowy Targ",    4, 3.1,  "28/10/2023 03:12", 2.97, "29/10/2023 13:48", 3.52, "28/10/2023 03:12", 3.61, "29/10/2023 13:38", 1.85, "28/10/2023 03:12", 2.03, "29/10/2023 13:48", "https://www.betexplorer.com/football/poland/iii-liga-group-iv/sokol-sieniawa-podhale-nowy-targ/ETmGMEhT/")
)

$startRow = 101
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $arr = New-Object 'object[,]' 1, 22
    for ($c = 0; $c -lt 22; $c++) {
        $arr[0, $c] = $rowData[$c]
    }
    $ws.Range("A${r}:V${r}").Value2 = $arr
}

Write-Host "Applied row reordering and appended rows 101-109."
